$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "Spotify Player - " / "embedded playlist from " / "Spotify." were three
#    separate runs with identical formatting; merge them into a single run.
#    Find/Replace re-writes the matched range as one run, picking up the
#    shared run formatting (sz/szCs) along the way.
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Spotify Player - embedded playlist from Spotify.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Spotify Player - embedded playlist from Spotify.", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) After "YouTube Player ..." there is an empty paragraph. Give it the text
#    "Share button - share music to social media.", then add two further new
#    paragraphs: "Live chat - section for listeners to speak to each other. "
#    and a blank one, followed by the pre-existing blank paragraph (so there
#    end up being two blank paragraphs before "Nonfunctional :").
#
#    To keep the new runs formatted like their neighbours (sz=24/szCs=24)
#    we append the new text directly onto the end of the preceding,
#    already-formatted paragraph (so it is typed into that run), then split
#    it off into its own paragraph with InsertParagraphAfter - mirroring how
#    a user typing in Word would inherit the preceding run's formatting.
# ---------------------------------------------------------------------------
$pYouTube = $d.Paragraphs.Item(6)
$pos = $pYouTube.Range.End - 1
$d.Range($pos, $pos).InsertAfter("Share button - share music to social media.") | Out-Null
$d.Range($pos, $pos).InsertParagraphAfter() | Out-Null

$pShare = $d.Paragraphs.Item(7)
$pos = $pShare.Range.End - 1
$d.Range($pos, $pos).InsertAfter("Live chat - section for listeners to speak to each other. ") | Out-Null
$d.Range($pos, $pos).InsertParagraphAfter() | Out-Null

# The paragraph that is now 9th is the original empty paragraph; duplicate it
# so there are two blank paragraphs before "Nonfunctional :".
$pBlank = $d.Paragraphs.Item(9)
$pBlank.Range.InsertParagraphAfter() | Out-Null

# ---------------------------------------------------------------------------
# 3) Before "Constraints:" there is an empty paragraph (after "Personal
#    information should be secure."). Give it the text "Terms and Services.",
#    then leave one blank paragraph after it (the original empty paragraph
#    ends up playing that role once the new text paragraph is split off).
# ---------------------------------------------------------------------------
$pPersonal = $d.Paragraphs.Item(15)
$pos = $pPersonal.Range.End - 1
$d.Range($pos, $pos).InsertAfter("Terms and Services.") | Out-Null
$d.Range($pos, $pos).InsertParagraphAfter() | Out-Null

# ---------------------------------------------------------------------------
# 4) "Resources - single person creating the website" gains a trailing period
#    and is split into two runs: "Resources - single person creating the "
#    and "website."
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "Resources " + [char]0x2013 + " single person creating the website",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Resources " + [char]0x2013 + " single person creating the website.", 2) | Out-Null

$pResources = $d.Paragraphs.Item($d.Paragraphs.Count)
$r = $pResources.Range
$r.End = $r.End - 1
$splitOffset = $r.Text.IndexOf("website.")
$splitPos = $r.Start + $splitOffset
$rAfter = $d.Range($splitPos, $r.End)
$rAfter.Bold = 1
$rAfter.Bold = 0

Write-Host "Edit complete."
